$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindingsTracker")

# Clear the "ITRCA Member that filed or reviewed(if BOT found) final evidence"
# column (M) for the data rows (2-17) per the updated Business Exception
# logging -- these no longer record an individual reviewer name.
$ws.Range("M2:M17").ClearContents()

# Update the view: scroll so column G is the left-most visible column and
# move the active selection to O7 (matches latest saved view state).
$ws.Activate()
$ws.Range("O7").Select()
$excel.ActiveWindow.ScrollColumn = 7
